$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-24 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-25 Thursday", 2) | Out-Null
$d.Content.Find.Execute("2+14=", $true, $false, $false, $false, $false, $true, 1, $false, "27+24=", 2) | Out-Null
$d.Content.Find.Execute("68-44=", $true, $false, $false, $false, $false, $true, 1, $false, "12+67=", 2) | Out-Null
$d.Content.Find.Execute("61+28=", $true, $false, $false, $false, $false, $true, 1, $false, "70-30=", 2) | Out-Null
$d.Content.Find.Execute("22+20=", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=", 2) | Out-Null
$d.Content.Find.Execute("60+26=", $true, $false, $false, $false, $false, $true, 1, $false, "31+40=", 2) | Out-Null
$d.Content.Find.Execute("65+6=", $true, $false, $false, $false, $false, $true, 1, $false, "72-33=", 2) | Out-Null
$d.Content.Find.Execute("15+23=", $true, $false, $false, $false, $false, $true, 1, $false, "1+48=", 2) | Out-Null
$d.Content.Find.Execute("9-3=", $true, $false, $false, $false, $false, $true, 1, $false, "76-8=", 2) | Out-Null
$d.Content.Find.Execute("31+12=", $true, $false, $false, $false, $false, $true, 1, $false, "37+8=", 2) | Out-Null
$d.Content.Find.Execute("18+25=", $true, $false, $false, $false, $false, $true, 1, $false, "93-42=", 2) | Out-Null
$d.Content.Find.Execute("26+34=", $true, $false, $false, $false, $false, $true, 1, $false, "78-19=", 2) | Out-Null
$d.Content.Find.Execute("94-38=", $true, $false, $false, $false, $false, $true, 1, $false, "73+7=", 2) | Out-Null
$d.Content.Find.Execute("34+21=", $true, $false, $false, $false, $false, $true, 1, $false, "86-11=", 2) | Out-Null
$d.Content.Find.Execute("8+56=", $true, $false, $false, $false, $false, $true, 1, $false, "61-28=", 2) | Out-Null
$d.Content.Find.Execute("76-31=", $true, $false, $false, $false, $false, $true, 1, $false, "1+62=", 2) | Out-Null
$d.Content.Find.Execute("98-30=", $true, $false, $false, $false, $false, $true, 1, $false, "14+79=", 2) | Out-Null
$d.Content.Find.Execute("65-48=", $true, $false, $false, $false, $false, $true, 1, $false, "39+50=", 2) | Out-Null
$d.Content.Find.Execute("18+67=", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=", 2) | Out-Null
$d.Content.Find.Execute("92-51=", $true, $false, $false, $false, $false, $true, 1, $false, "54-23=", 2) | Out-Null
$d.Content.Find.Execute("70-3=", $true, $false, $false, $false, $false, $true, 1, $false, "40+6=", 2) | Out-Null
$d.Content.Find.Execute("44+20=", $true, $false, $false, $false, $false, $true, 1, $false, "74+18=", 2) | Out-Null
$d.Content.Find.Execute("13-2=", $true, $false, $false, $false, $false, $true, 1, $false, "13-13=", 2) | Out-Null
$d.Content.Find.Execute("66-28=", $true, $false, $false, $false, $false, $true, 1, $false, "77-22=", 2) | Out-Null
$d.Content.Find.Execute("37-0=", $true, $false, $false, $false, $false, $true, 1, $false, "68-0=", 2) | Out-Null
$d.Content.Find.Execute("53-4=", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=", 2) | Out-Null
$d.Content.Find.Execute("39+25=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2) | Out-Null
$d.Content.Find.Execute("14+22=", $true, $false, $false, $false, $false, $true, 1, $false, "57+16=", 2) | Out-Null
$d.Content.Find.Execute("70-43=", $true, $false, $false, $false, $false, $true, 1, $false, "93-60=", 2) | Out-Null
$d.Content.Find.Execute("97-24=", $true, $false, $false, $false, $false, $true, 1, $false, "93-84=", 2) | Out-Null
$d.Content.Find.Execute("18+55=", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=", 2) | Out-Null
$d.Content.Find.Execute("79-44=", $true, $false, $false, $false, $false, $true, 1, $false, "94-91=", 2) | Out-Null
$d.Content.Find.Execute("66-25=", $true, $false, $false, $false, $false, $true, 1, $false, "37+13=", 2) | Out-Null
$d.Content.Find.Execute("6+81=", $true, $false, $false, $false, $false, $true, 1, $false, "84-41=", 2) | Out-Null
$d.Content.Find.Execute("70-40=", $true, $false, $false, $false, $false, $true, 1, $false, "72+6=", 2) | Out-Null
$d.Content.Find.Execute("10+70=", $true, $false, $false, $false, $false, $true, 1, $false, "55+13=", 2) | Out-Null
$d.Content.Find.Execute("36+56=", $true, $false, $false, $false, $false, $true, 1, $false, "0+53=", 2) | Out-Null
$d.Content.Find.Execute("65-44=", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=", 2) | Out-Null
$d.Content.Find.Execute("64-56=", $true, $false, $false, $false, $false, $true, 1, $false, "72-14=", 2) | Out-Null
$d.Content.Find.Execute("9+79=", $true, $false, $false, $false, $false, $true, 1, $false, "19+2=", 2) | Out-Null
$d.Content.Find.Execute("56+29=", $true, $false, $false, $false, $false, $true, 1, $false, "81-11=", 2) | Out-Null
$d.Content.Find.Execute("27+18=", $true, $false, $false, $false, $false, $true, 1, $false, "73+23=", 2) | Out-Null
$d.Content.Find.Execute("58+14=", $true, $false, $false, $false, $false, $true, 1, $false, "9+1=", 2) | Out-Null
$d.Content.Find.Execute("58+40=", $true, $false, $false, $false, $false, $true, 1, $false, "11+3=", 2) | Out-Null
$d.Content.Find.Execute("39-34=", $true, $false, $false, $false, $false, $true, 1, $false, "5+92=", 2) | Out-Null
$d.Content.Find.Execute("45+52=", $true, $false, $false, $false, $false, $true, 1, $false, "93-33=", 2) | Out-Null
$d.Content.Find.Execute("4+7=", $true, $false, $false, $false, $false, $true, 1, $false, "85+1=", 2) | Out-Null
$d.Content.Find.Execute("55-23=", $true, $false, $false, $false, $false, $true, 1, $false, "1+97=", 2) | Out-Null
$d.Content.Find.Execute("14+31=", $true, $false, $false, $false, $false, $true, 1, $false, "95-89=", 2) | Out-Null
$d.Content.Find.Execute("54+22=", $true, $false, $false, $false, $false, $true, 1, $false, "41-13=", 2) | Out-Null
$d.Content.Find.Execute("62+16=", $true, $false, $false, $false, $false, $true, 1, $false, "58+1=", 2) | Out-Null
$d.Content.Find.Execute("68+12=", $true, $false, $false, $false, $false, $true, 1, $false, "93-67=", 2) | Out-Null
$d.Content.Find.Execute("2+76=", $true, $false, $false, $false, $false, $true, 1, $false, "82+7=", 2) | Out-Null
$d.Content.Find.Execute("82-2=", $true, $false, $false, $false, $false, $true, 1, $false, "69-12=", 2) | Out-Null
$d.Content.Find.Execute("66+18=", $true, $false, $false, $false, $false, $true, 1, $false, "50-25=", 2) | Out-Null
$d.Content.Find.Execute("70+29=", $true, $false, $false, $false, $false, $true, 1, $false, "63-39=", 2) | Out-Null
$d.Content.Find.Execute("26+12=", $true, $false, $false, $false, $false, $true, 1, $false, "66-45=", 2) | Out-Null
$d.Content.Find.Execute("87-58=", $true, $false, $false, $false, $false, $true, 1, $false, "62-6=", 2) | Out-Null
$d.Content.Find.Execute("68-58=", $true, $false, $false, $false, $false, $true, 1, $false, "48+10=", 2) | Out-Null
$d.Content.Find.Execute("31+24=", $true, $false, $false, $false, $false, $true, 1, $false, "36+51=", 2) | Out-Null
$d.Content.Find.Execute("5+19=", $true, $false, $false, $false, $false, $true, 1, $false, "18+13=", 2) | Out-Null
$d.Content.Find.Execute("51-28=", $true, $false, $false, $false, $false, $true, 1, $false, "58-38=", 2) | Out-Null
$d.Content.Find.Execute("31+49=", $true, $false, $false, $false, $false, $true, 1, $false, "84-79=", 2) | Out-Null
$d.Content.Find.Execute("58-4=", $true, $false, $false, $false, $false, $true, 1, $false, "13+27=", 2) | Out-Null
$d.Content.Find.Execute("36+20=", $true, $false, $false, $false, $false, $true, 1, $false, "50+26=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $true, $false, $false, $false, $false, $true, 1, $false, "98-65=", 2) | Out-Null
$d.Content.Find.Execute("38+7=", $true, $false, $false, $false, $false, $true, 1, $false, "53+37=", 2) | Out-Null
$d.Content.Find.Execute("61-18=", $true, $false, $false, $false, $false, $true, 1, $false, "94-20=", 2) | Out-Null
$d.Content.Find.Execute("79-4=", $true, $false, $false, $false, $false, $true, 1, $false, "60-41=", 2) | Out-Null
$d.Content.Find.Execute("37+47=", $true, $false, $false, $false, $false, $true, 1, $false, "33+40=", 2) | Out-Null
$d.Content.Find.Execute("7+37=", $true, $false, $false, $false, $false, $true, 1, $false, "83-75=", 2) | Out-Null
$d.Content.Find.Execute("67-63=", $true, $false, $false, $false, $false, $true, 1, $false, "45+22=", 2) | Out-Null
$d.Content.Find.Execute("93+4=", $true, $false, $false, $false, $false, $true, 1, $false, "82-33=", 2) | Out-Null
$d.Content.Find.Execute("7+4=", $true, $false, $false, $false, $false, $true, 1, $false, "52-47=", 2) | Out-Null
$d.Content.Find.Execute("53-51=", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("0+38=", $true, $false, $false, $false, $false, $true, 1, $false, "47-19=", 2) | Out-Null
$d.Content.Find.Execute("9-1=", $true, $false, $false, $false, $false, $true, 1, $false, "63+24=", 2) | Out-Null
$d.Content.Find.Execute("99-93=", $true, $false, $false, $false, $false, $true, 1, $false, "4+94=", 2) | Out-Null
$d.Content.Find.Execute("63-34=", $true, $false, $false, $false, $false, $true, 1, $false, "21+6=", 2) | Out-Null
$d.Content.Find.Execute("74-16=", $true, $false, $false, $false, $false, $true, 1, $false, "88-37=", 2) | Out-Null
$d.Content.Find.Execute("18+27=", $true, $false, $false, $false, $false, $true, 1, $false, "12+78=", 2) | Out-Null
$d.Content.Find.Execute("39+59=", $true, $false, $false, $false, $false, $true, 1, $false, "58+5=", 2) | Out-Null
$d.Content.Find.Execute("37+25=", $true, $false, $false, $false, $false, $true, 1, $false, "85-16=", 2) | Out-Null
$d.Content.Find.Execute("40+23=", $true, $false, $false, $false, $false, $true, 1, $false, "70+15=", 2) | Out-Null
$d.Content.Find.Execute("53-2=", $true, $false, $false, $false, $false, $true, 1, $false, "13+13=", 2) | Out-Null
$d.Content.Find.Execute("62-29=", $true, $false, $false, $false, $false, $true, 1, $false, "65-62=", 2) | Out-Null
$d.Content.Find.Execute("73-8=", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=", 2) | Out-Null
$d.Content.Find.Execute("60-25=", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=", 2) | Out-Null
$d.Content.Find.Execute("0+1=", $true, $false, $false, $false, $false, $true, 1, $false, "4+57=", 2) | Out-Null
$d.Content.Find.Execute("64+17=", $true, $false, $false, $false, $false, $true, 1, $false, "67-20=", 2) | Out-Null
$d.Content.Find.Execute("10+59=", $true, $false, $false, $false, $false, $true, 1, $false, "34-13=", 2) | Out-Null
$d.Content.Find.Execute("57-31=", $true, $false, $false, $false, $false, $true, 1, $false, "78-4=", 2) | Out-Null
$d.Content.Find.Execute("38-0=", $true, $false, $false, $false, $false, $true, 1, $false, "24+38=", 2) | Out-Null
$d.Content.Find.Execute("42+11=", $true, $false, $false, $false, $false, $true, 1, $false, "59+35=", 2) | Out-Null
$d.Content.Find.Execute("52-28=", $true, $false, $false, $false, $false, $true, 1, $false, "73-47=", 2) | Out-Null
$d.Content.Find.Execute("2+63=", $true, $false, $false, $false, $false, $true, 1, $false, "95-40=", 2) | Out-Null
$d.Content.Find.Execute("77+21=", $true, $false, $false, $false, $false, $true, 1, $false, "70+23=", 2) | Out-Null
$d.Content.Find.Execute("27-17=", $true, $false, $false, $false, $false, $true, 1, $false, "3+18=", 2) | Out-Null
$d.Content.Find.Execute("76-45=", $true, $false, $false, $false, $false, $true, 1, $false, "41+28=", 2) | Out-Null
$d.Content.Find.Execute("76-38=", $true, $false, $false, $false, $false, $true, 1, $false, "18-14=", 2) | Out-Null
$d.Content.Find.Execute("77-41=", $true, $false, $false, $false, $false, $true, 1, $false, "99-58=", 2) | Out-Null
